$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "alpha2F"

# Update precise floating point values in row 13
$ws.Range("D13").Value = 0.9950778694753355
$ws.Range("H13").Value = 0.9959845843276111
$ws.Range("J13").Value = 0.9950778694753355
$ws.Range("M13").Value = 0.9954104765740497
